# Push - 21.11 - Removed all folder paths from Config
#
# The Orchestrator folder ("UatRPA/BAW/BA Holidays/BA Holidays") used to be
# hard-coded into the Config workbook in two places:
#   - Settings!B3  (OrchestratorQueueFolder)
#   - Assets!C2:C35 (OrchestratorAssetFolder column, one per asset row)
# Both are cleared out below so the workbook no longer references any
# specific folder path.

$wb = $excel.ActiveWorkbook

# --- Settings sheet: clear the OrchestratorQueueFolder value --------------
$wsSettings = $wb.Worksheets.Item("Settings")
[void]$wsSettings.Activate()
$wsSettings.Range("B3").ClearContents()
$wsSettings.Rows.Item(3).RowHeight = 45
$wsSettings.Rows.Item(5).RowHeight = 30
[void]$wsSettings.Range("B7").Select()

# --- Constants sheet: no value changes, just row-height touch-ups ---------
$wsConstants = $wb.Worksheets.Item("Constants")
[void]$wsConstants.Activate()
$wsConstants.Rows.Item(2).RowHeight = 30
$wsConstants.Rows.Item(3).RowHeight = 45
$wsConstants.Rows.Item(17).RowHeight = 45
[void]$wsConstants.Range("C25").Select()

# --- Assets sheet: clear the OrchestratorAssetFolder column ---------------
$wsAssets = $wb.Worksheets.Item("Assets")
[void]$wsAssets.Activate()
$wsAssets.Range("C2:C35").ClearContents()
[void]$wsAssets.Range("C14").Select()
